$wb = $excel.ActiveWorkbook

# The new convention/event row that was scraped and appended to the
# "展览" (Exhibition) and "全部类型" (All types) sheets.
$targetSheets = @("展览", "全部类型")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # A2 mirrors A1's numeric index style (bold, bordered, centered),
    # so copy A1's format down before writing the new index value.
    $ws.Cells.Item(1, 1).Copy()
    $ws.Cells.Item(2, 1).PasteSpecial(-4122)
    $ws.Cells.Item(2, 1).Value = 1

    # B2 is a literal date-like text string ("2024-03-30"), not a real
    # Excel date, so force text interpretation with a leading apostrophe
    # to avoid automatic date conversion.
    $ws.Cells.Item(2, 2).Value = "'2024-03-30"

    $ws.Cells.Item(2, 3).Value = "丽水·2024首届TCT国风动漫游戏嘉年华"
    $ws.Cells.Item(2, 4).Value = "城北街11号华东药用植物园百药谷 中医药文化展览馆"
    $ws.Cells.Item(2, 5).Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Cells.Item(2, 6).Value = 1
    $ws.Cells.Item(2, 7).Value = 68
    $ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82254"
    $ws.Cells.Item(2, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/ZPmB9ko51709521574642.jpeg"
}
